# Generate Report for Handoff
#
# Refreshes the localization-status report after a fresh handoff run:
#  - zh-cn / de-de sheets: rows 4-7 (the items that were still "low"
#    priority / not yet localized for this language) flip to "ht" and
#    pick up the new handoff timestamp.
#  - Overview sheet: the "Latest HO Xliff Generate Date" column mirrors
#    the de-de handoff timestamp for those same rows.

$wb = $excel.ActiveWorkbook

$overview = $wb.Worksheets.Item("Overview")
$zhcn     = $wb.Worksheets.Item("zh-cn")
$dede     = $wb.Worksheets.Item("de-de")

# zh-cn: Priority (E) low -> ht, Latest Handoff Datetime (H) refreshed
$zhcn.Range("E4:E7").Value = "ht"
$zhcn.Range("H4:H7").Value = "2016-08-24 02:31:28"

# de-de: Priority (E) low -> ht, Latest Handoff Datetime (H) refreshed
$dede.Range("E4:E7").Value = "ht"
$dede.Range("H4:H7").Value = "2016-08-24 02:31:33"

# Overview: Latest HO Xliff Generate Date (G) refreshed to match
$overview.Range("G4:G7").Value = "2016-08-24 02:31:33"
